# Tabs with H-Score results include the scored marker name, in support of #49
$wb = $excel.ActiveWorkbook

# Rename the H-Score sheets to include the scored marker name (PDL1)
$wsHScore = $wb.Worksheets.Item("H-Score")
$wsHScore.Name = "H-Score, PDL1"

$wsHScoreCD8 = $wb.Worksheets.Item("H-Score CD8+")
$wsHScoreCD8.Name = "H-Score, PDL1, CD8+"

$wsHScoreCK = $wb.Worksheets.Item("H-Score CK+_Membrane PDL1 (Opal")
$wsHScoreCK.Name = "H-Score, PDL1, CK+_PDL1>1"

# Update the Print_Titles defined names to use the new sheet names
foreach ($n in $wb.Names) {
    if ($n.Name -eq "H-Score, PDL1!Print_Titles") {
        $n.RefersTo = "='H-Score, PDL1'!`$1:`$3"
    } elseif ($n.Name -eq "H-Score, PDL1, CD8+!Print_Titles") {
        $n.RefersTo = "='H-Score, PDL1, CD8+'!`$1:`$3"
    } elseif ($n.Name -eq "H-Score, PDL1, CK+_PDL1>1!Print_Titles") {
        $n.RefersTo = "='H-Score, PDL1, CK+_PDL1>1'!`$1:`$3"
    }
}

# Update sheet selections to match the saved view state. Selecting a range on
# a sheet also makes that sheet active, so the sheet that should remain the
# active/selected tab at the end is handled last.
$wsMeanExpr = $wb.Worksheets.Item("Mean Expression")
$wsMeanExpr.Range("G2").Select() | Out-Null

$wsHScore.Range("A1:A3").Select() | Out-Null

$wsHScoreCD8.Range("A1:A3").Select() | Out-Null

$wsNearest = $wb.Worksheets.Item("Nearest Neighbors")
$wsNearest.Range("A1:A2").Select() | Out-Null

# This sheet ends up as the active / tabSelected sheet, so select it last.
$wsHScoreCK.Range("A1:A3").Select() | Out-Null
